$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write new shared-string values in the exact order they must
# first appear, so the shared-strings table receives them in the order
# required by the target file (indices 40..48).
#
# A couple of the new date-like strings ("11.3.21" and "4.3.21") would be
# auto-parsed by the Value setter into date serial numbers (since the
# first numeric component is <=12 and so looks like a valid month). To
# keep them as plain text (matching the original file's convention of
# storing dates as text), we build them via a formula that evaluates to
# the literal text, then convert that formula to a static value in place
# via PasteSpecial(values). This avoids the date auto-detection entirely
# and does not touch styles.xml.

# 1) "Y" -> J3
$ws.Range("J3").Value = "Y"

# 2) "11.3.21" -> A20  (needs the formula/paste-values trick)
$ws.Range("A20").Formula = "=""11.3.21"""
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4163) | Out-Null

# 3) "Full run on myself to check timing" -> U21
$ws.Range("U21").Value = "Full run on myself to check timing"

# 4) "14.3.21" -> A21 (first component 14 > 12, stays text naturally)
$ws.Range("A21").Value = "14.3.21"

# 5) "Full run on myself to check manual tests and timing" -> U20
$ws.Range("U20").Value = "Full run on myself to check manual tests and timing"

# 6) "4.3.21" -> A3 (needs the formula/paste-values trick)
$ws.Range("A3").Formula = "=""4.3.21"""
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4163) | Out-Null

# 7) "right" -> F3
$ws.Range("F3").Value = "right"

# 8) "glasses" -> G3
$ws.Range("G3").Value = "glasses"

# 9) "F" -> I3
$ws.Range("I3").Value = "F"

# --- Step 2: remaining new cells that reuse already-existing shared
# strings, or that hold plain numbers. Order does not matter for these.
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 27
$ws.Range("L3").Value = "N"
$ws.Range("M3").Value = "N"
$ws.Range("N3").Value = "N"

$ws.Range("B20").Value = 1007
$ws.Range("P20").Value = "118B"
$ws.Range("T20").Value = "Khen heller"

$ws.Range("B21").Value = 1008
$ws.Range("P21").Value = "118B"
$ws.Range("T21").Value = "Khen heller"

# --- Step 3: match row heights used throughout the rest of the sheet.
$ws.Rows.Item(20).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75

# --- Step 4: update the selected cell shown when the workbook is opened.
$ws.Range("C4").Select() | Out-Null
